# Iraq League.xlsx update (11-06-2024 21:19)
# The source feed re-ordered a handful of fixtures that share the same
# date/matchday block, which shows up as each pair of adjacent data rows
# having their match details (id/home-away teams/odds/etc., columns B:AD)
# swapped while the running row number in column A stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is a pair of worksheet rows whose B:AD contents (everything
# except the column-A row counter) need to be swapped with one another.
$rowPairs = @(
    @(41, 42),
    @(58, 59),
    @(69, 70),
    @(78, 79),
    @(89, 90),
    @(219, 220),
    @(221, 222),
    @(223, 224)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}
